# Update "想去人数" (interested count) values in column F across the
# "展览" (Exhibitions), "本地生活" (Local Life), and "全部类型" (All Types)
# worksheets, reflecting the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 57
$ws1.Range("F4").Value = 2012
$ws1.Range("F5").Value = 5883
$ws1.Range("F8").Value = 3334
$ws1.Range("F11").Value = 1396
$ws1.Range("F12").Value = 4643
$ws1.Range("F13").Value = 1100
$ws1.Range("F14").Value = 1749
$ws1.Range("F15").Value = 11
$ws1.Range("F18").Value = 203
$ws1.Range("F20").Value = 1042
$ws1.Range("F21").Value = 316
$ws1.Range("F28").Value = 1135
$ws1.Range("F30").Value = 104
$ws1.Range("F31").Value = 221
$ws1.Range("F32").Value = 451
$ws1.Range("F35").Value = 1784
$ws1.Range("F37").Value = 1067
$ws1.Range("F42").Value = 410
$ws1.Range("F47").Value = 439

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 800

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 800
$ws4.Range("F4").Value = 57
$ws4.Range("F5").Value = 2012
$ws4.Range("F6").Value = 5883
$ws4.Range("F9").Value = 3334
$ws4.Range("F11").Value = 1396
$ws4.Range("F12").Value = 4643
$ws4.Range("F13").Value = 1749
$ws4.Range("F14").Value = 11
$ws4.Range("F20").Value = 203
$ws4.Range("F23").Value = 1042
$ws4.Range("F24").Value = 316
$ws4.Range("F28").Value = 1135
$ws4.Range("F30").Value = 104
$ws4.Range("F31").Value = 221
$ws4.Range("F33").Value = 1784
$ws4.Range("F35").Value = 1067
$ws4.Range("F41").Value = 410
$ws4.Range("F44").Value = 439

$wb.Save()
